$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "COBERTURA EM DIAS" column (old AP) entirely; columns after it shift left
$ws.Columns("AP:AP").Delete()

# Update header row (row 1) text: de-accent/tighten "%" labels, fix casing, drop stray spaces
$ws.Range("A1").Value = "DATA_HORA_ANALISE"
$ws.Range("B1").Value = "TOTAL DE SKUs"
$ws.Range("C1").Value = "TOTAL SKU COM VENDA ACIMA DE 1 ANO"
$ws.Range("D1").Value = "%SKU COM VENDA ACIMA DE 1 ANO"
$ws.Range("E1").Value = "TOTAL SKU COM VENDA SOMENTE NO ULTIMO ANO"
$ws.Range("F1").Value = "%SKU COM VENDA SOMENTE NO ULTIMO ANO"
$ws.Range("G1").Value = "TOTAL SKU COM ESTOQUE ZERO"
$ws.Range("H1").Value = "%SKU COM ESTOQUE ZERO"
$ws.Range("I1").Value = "TOTAL SKU COM ESTOQUE POSITIVO"
$ws.Range("J1").Value = "%SKU COM ESTOQUE POSITIVO"
$ws.Range("K1").Value = "CUSTO TOTAL ESTOQUE POSITIVO"
$ws.Range("L1").Value = "TOTAL SKU COM ESTOQUE NEGATIVO"
$ws.Range("M1").Value = "%SKU COM ESTOQUE NEGATIVO"
$ws.Range("N1").Value = "CUSTO TOTAL ESTOQUE NEGATIVO"
$ws.Range("O1").Value = "TOTAL SKU INATIVO COM SALDO"
$ws.Range("P1").Value = "%SKU INATIVO COM SALDO"
$ws.Range("Q1").Value = "CUSTO TOTAL INATIVO COM SALDO"
$ws.Range("R1").Value = "TOTAL SKU INATIVO SEM SALDO"
$ws.Range("S1").Value = "%SKU INATIVO SEM SALDO"
$ws.Range("T1").Value = "TOTAL SKU ATIVO COM SALDO"
$ws.Range("U1").Value = "%SKU ATIVO COM SALDO"
$ws.Range("V1").Value = "CUSTO TOTAL ATIVO COM SALDO"
$ws.Range("W1").Value = "TOTAL SKU ATIVO SEM SALDO"
$ws.Range("X1").Value = "%SKU ATIVO SEM SALDO"
$ws.Range("Y1").Value = "TOTAL SKU SEM VENDA COM SALDO"
$ws.Range("Z1").Value = "%SKU SEM VENDA COM SALDO"
$ws.Range("AA1").Value = "CUSTO TOTAL SEM VENDA COM SALDO"
$ws.Range("AB1").Value = "TOTAL SKU SEM VENDA SEM SALDO"
$ws.Range("AC1").Value = "%SKU SEM VENDA SEM SALDO"
$ws.Range("AD1").Value = "TOTAL SKU GRUPO A"
$ws.Range("AE1").Value = "TOTAL SKU GRUPO B"
$ws.Range("AF1").Value = "TOTAL SKU GRUPO C"
$ws.Range("AG1").Value = "%SKU GRUPO A"
$ws.Range("AH1").Value = "%SKU GRUPO B"
$ws.Range("AI1").Value = "%SKU GRUPO C"
$ws.Range("AJ1").Value = "TOTAL VENDA GRUPO A"
$ws.Range("AK1").Value = "TOTAL VENDA GRUPO B"
$ws.Range("AL1").Value = "TOTAL VENDA GRUPO C"
$ws.Range("AM1").Value = "%VENDA GRUPO A"
$ws.Range("AN1").Value = "%VENDA GRUPO B"
$ws.Range("AO1").Value = "%VENDA GRUPO C"
$ws.Range("AP1").Value = "COBERTURA EM DIAS GRUPO A"
$ws.Range("AQ1").Value = "COBERTURA EM DIAS GRUPO B"
$ws.Range("AR1").Value = "COBERTURA EM DIAS GRUPO C"

# Refresh data row (row 2) with the latest snapshot
$ws.Range("A2").Value = "2025-05-09 18:54:59"
$ws.Range("B2").Value = 16142
$ws.Range("C2").Value = 11643
$ws.Range("D2").Value = 72.12860859868667
$ws.Range("E2").Value = 2277
$ws.Range("F2").Value = 14.10605872878206
$ws.Range("G2").Value = 3094
$ws.Range("H2").Value = 19.1673894189072
$ws.Range("I2").Value = 9526
$ws.Range("J2").Value = 59.01375294263412
$ws.Range("K2").Value = 3036032.85
$ws.Range("L2").Value = 3522
$ws.Range("M2").Value = 21.81885763845868
$ws.Range("N2").Value = 1072804.73
$ws.Range("O2").Value = 4623
$ws.Range("P2").Value = 28.63957378267873
$ws.Range("Q2").Value = 473467.8999999999
$ws.Range("R2").Value = 3502
$ws.Range("S2").Value = 21.69495725436749
$ws.Range("T2").Value = 3484
$ws.Range("U2").Value = 21.58344690868542
$ws.Range("V2").Value = 2422321.85
$ws.Range("W2").Value = 2311
$ws.Range("X2").Value = 14.31668938173708
$ws.Range("Y2").Value = 1419
$ws.Range("Z2").Value = 8.790732251269979
$ws.Range("AA2").Value = 140243.1
$ws.Range("AB2").Value = 803
$ws.Range("AC2").Value = 4.974600421261306
$ws.Range("AD2").Value = 463
$ws.Range("AE2").Value = 878
$ws.Range("AF2").Value = 1389
$ws.Range("AG2").Value = 16.95970695970696
$ws.Range("AH2").Value = 32.16117216117216
$ws.Range("AI2").Value = 50.87912087912088
$ws.Range("AJ2").Value = 1422342.64
$ws.Range("AK2").Value = 266803.27
$ws.Range("AL2").Value = 89026.19
$ws.Range("AM2").Value = 79.98903143289674
$ws.Range("AN2").Value = 15.00435587758913
$ws.Range("AO2").Value = 5.006612689514137
$ws.Range("AP2").Value = 47.69122524876872
$ws.Range("AQ2").Value = 224.0439682246444
$ws.Range("AR2").Value = 568.357470858354
